$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for data rows 2-43
# from 2025-04-14 (45761) to 2025-04-15 (45762).
for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45761) {
        $cell.Value2 = 45762
    }
}
